$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.880.43'
$ws.Range("E2").Value = "'  -1.88%  "
$ws.Range("D3").Value = '3.406.16'
$ws.Range("E3").Value = "'  -1.43%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'578.81"
$ws.Range("E5").Value = "'  -2.34%  "
$ws.Range("D6").Value = "'171.18"
$ws.Range("E6").Value = "'  -4.58%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "'  -3.97%  "
$ws.Range("D9").Value = '3.405.41'
$ws.Range("E9").Value = "'  -1.41%  "
$ws.Range("E10").Value = "'  -7.03%  "
$ws.Range("D11").Value = "'6.82"
$ws.Range("E11").Value = "'  -1.92%  "
$ws.Range("D12").Value = "'0.404"
$ws.Range("E12").Value = "'  -5.69%  "
$ws.Range("D13").Value = '4.002.79'
$ws.Range("E13").Value = "'  -1.28%  "
$ws.Range("E14").Value = "'  -0.72%  "
$ws.Range("D15").Value = "'29.55"
$ws.Range("E15").Value = "'  -7.56%  "
$ws.Range("D16").Value = '65.946.54'
$ws.Range("E16").Value = "'  -1.71%  "
$ws.Range("D17").Value = "'0.0000168"
$ws.Range("E17").Value = "'  -4.66%  "
$ws.Range("D18").Value = '3.406.77'
$ws.Range("E18").Value = "'  -1.35%  "
$ws.Range("E19").Value = "'  -5.64%  "
$ws.Range("D20").Value = "'13.56"
$ws.Range("E20").Value = "'  -3.77%  "
$ws.Range("D21").Value = "'361.33"
$ws.Range("E21").Value = "'  -7.82%  "
$ws.Range("D22").Value = "'7.59"
$ws.Range("E22").Value = "'  -3.80%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("D24").Value = "'5.71"
$ws.Range("E24").Value = "'  -1.08%  "
$ws.Range("D25").Value = "'70.86"
$ws.Range("E25").Value = "'  -0.97%  "
$ws.Range("D26").Value = "'0.522"
$ws.Range("E26").Value = "'  -2.63%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "'  -2.92%  "
$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = "'  -7.60%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  -0.09%  "
$ws.Range("D31").Value = "'23.59"
$ws.Range("E31").Value = "'  +0.68%  "
$ws.Range("D32").Value = "'5.72"
$ws.Range("D33").Value = "'1.96"
$ws.Range("E33").Value = "'  -3.99%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("E35").Value = "'  -8.30%  "
$ws.Range("D36").Value = "'6.95"
$ws.Range("E36").Value = "'  -4.79%  "
$ws.Range("E37").Value = "'  -3.55%  "
$ws.Range("D38").Value = "'159.97"
$ws.Range("E38").Value = "'  -0.58%  "
$ws.Range("D39").Value = "'28.85"
$ws.Range("E39").Value = "'  +10.94%  "
$ws.Range("D40").Value = "'0.875"
$ws.Range("E40").Value = "'  -0.14%  "
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "'  -7.23%  "
$ws.Range("E42").Value = "'  -7.26%  "
$ws.Range("D43").Value = '2.680.30'
$ws.Range("E43").Value = "'  -2.64%  "
$ws.Range("D44").Value = "'4.35"
$ws.Range("E44").Value = "'  -6.14%  "
$ws.Range("D45").Value = "'6.29"
$ws.Range("E45").Value = "'  -6.66%  "
$ws.Range("D46").Value = "'0.0672"
$ws.Range("E46").Value = "'  -6.18%  "
$ws.Range("D47").Value = "'39.85"
$ws.Range("E47").Value = "'  -3.56%  "
$ws.Range("D48").Value = "'0.0285"
$ws.Range("E48").Value = "'  -4.02%  "
$ws.Range("D49").Value = "'23.73"
$ws.Range("E49").Value = "'  -9.60%  "
$ws.Range("D50").Value = "'303.06"
$ws.Range("E50").Value = "'  -6.43%  "
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = "'0.807"
$ws.Range("E51").Value = "'  -4.05%  "
